$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.650.82"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.153.90"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'531.40"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'139.79"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +15.28%  "
$ws.Range("D9").Value = "'7.32"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").Value = "3.696.08"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "'25.86"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "58.685.72"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").Value = "3.163.94"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'12.99"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "'8.14"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'371.58"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'0.526"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").Value = "'69.65"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'1.01"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'8.26"
$ws.Range("E28").Value = "  +12.84%  "
$ws.Range("D29").Value = "0.0₃0866"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.13"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.88"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.11"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'158.69"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").Value = "2.634.88"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").Value = "'4.24"
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("D43").Value = "'39.02"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "3.194.05"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +13.09%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "'0.979"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "'20.23"
$ws.Range("E51").Value = "  +1.86%  "
